$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $text) {
    $escaped = $text -replace '"', '""'
    $ws.Range("ZZ1").Formula = '="' + $escaped + '"'
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellAddr).PasteSpecial(-4163)
}

Set-TextValue "E2" "2026-02-10 07:48:25"
Set-TextValue "E3" "2026-02-10 07:48:27"
Set-TextValue "G3" "187 cm"
Set-TextValue "I3" "9.1 mm"
Set-TextValue "M3" "-0.1 °C 7:29 TU"
Set-TextValue "E4" "2026-02-10 07:48:30"
Set-TextValue "H4" "85%"
Set-TextValue "K4" "0.0 MJ/m2"
Set-TextValue "O4" "8.8 °C"
Set-TextValue "E5" "2026-02-10 07:48:32"
Set-TextValue "G5" "133 cm"
Set-TextValue "M5" "0.1 °C 7:28 TU"
Set-TextValue "E6" "2026-02-10 07:48:34"
Set-TextValue "O6" "7.1 °C"
Set-TextValue "E7" "2026-02-10 07:48:37"
Set-TextValue "E8" "2026-02-10 07:48:39"
Set-TextValue "K8" "0.1 MJ/m2"
Set-TextValue "E9" "2026-02-10 07:48:42"
Set-TextValue "E10" "2026-02-10 07:48:44"
Set-TextValue "N10" "4.5 °C 7:00 TU"
Set-TextValue "O10" "6.6 °C"
Set-TextValue "E11" "2026-02-10 07:48:46"
Set-TextValue "E12" "2026-02-10 07:48:49"
Set-TextValue "N12" "4.9 °C 7:03 TU"
Set-TextValue "E13" "2026-02-10 07:48:51"
Set-TextValue "I13" "2.6 mm"
Set-TextValue "J13" "1008.5 hPa"
Set-TextValue "N13" "2.3 °C 7:02 TU"
Set-TextValue "E14" "2026-02-10 07:48:54"
Set-TextValue "E15" "2026-02-10 07:48:56"
Set-TextValue "I15" "0.2 mm"
Set-TextValue "E16" "2026-02-10 07:48:58"
Set-TextValue "H16" "88%"
Set-TextValue "I16" "14.0 mm"
Set-TextValue "K16" "0.1 MJ/m2"
Set-TextValue "M16" "0.6 °C 7:00 TU"
Set-TextValue "E17" "2026-02-10 07:49:01"
Set-TextValue "K17" "0.1 MJ/m2"
Set-TextValue "O17" "3.1 °C"
Set-TextValue "E18" "2026-02-10 07:49:03"
Set-TextValue "J18" "1005.3 hPa"
Set-TextValue "O18" "6.8 °C"
Set-TextValue "E19" "2026-02-10 07:49:06"
Set-TextValue "E20" "2026-02-10 07:49:08"
Set-TextValue "M20" "-0.2 °C 7:27 TU"
Set-TextValue "E21" "2026-02-10 07:49:10"
Set-TextValue "I21" "3.9 mm"
Set-TextValue "E22" "2026-02-10 07:49:13"
Set-TextValue "M22" "-0.7 °C 7:29 TU"
Set-TextValue "O22" "-1.8 °C"
Set-TextValue "E23" "2026-02-10 07:49:15"
Set-TextValue "M23" "1.1 °C 7:29 TU"
Set-TextValue "O23" "-0.3 °C"
Set-TextValue "E24" "2026-02-10 07:49:18"
Set-TextValue "N24" "8.0 °C 7:28 TU"
Set-TextValue "E25" "2026-02-10 07:49:20"
Set-TextValue "H25" "94%"
Set-TextValue "O25" "-0.5 °C"
Set-TextValue "E26" "2026-02-10 07:49:23"
Set-TextValue "H26" "85%"
Set-TextValue "O26" "3.5 °C"
Set-TextValue "E27" "2026-02-10 07:49:25"
Set-TextValue "G27" "174 cm"
Set-TextValue "M27" "0.3 °C 7:15 TU"
Set-TextValue "E28" "2026-02-10 07:49:28"
Set-TextValue "O28" "4.7 °C"
Set-TextValue "E29" "2026-02-10 07:49:30"
Set-TextValue "O29" "8.2 °C"
Set-TextValue "E30" "2026-02-10 07:49:33"
Set-TextValue "O30" "7.2 °C"
Set-TextValue "E31" "2026-02-10 07:49:35"
Set-TextValue "K31" "0.0 MJ/m2"
Set-TextValue "E32" "2026-02-10 07:49:38"
Set-TextValue "I32" "0.7 mm"
Set-TextValue "M32" "8.6 °C 7:17 TU"
Set-TextValue "O32" "7.8 °C"
Set-TextValue "E33" "2026-02-10 07:49:40"
Set-TextValue "I33" "6.3 mm"
Set-TextValue "J33" "1008.0 hPa"
Set-TextValue "E34" "2026-02-10 07:49:43"
Set-TextValue "H34" "82%"
Set-TextValue "I34" "3.3 mm"
Set-TextValue "K34" "0.1 MJ/m2"
Set-TextValue "E35" "2026-02-10 07:49:45"
Set-TextValue "M35" "11.7 °C 7:19 TU"
Set-TextValue "O35" "10.5 °C"
Set-TextValue "E36" "2026-02-10 07:49:48"
Set-TextValue "J36" "1005.1 hPa"
Set-TextValue "N36" "5.8 °C 7:12 TU"
Set-TextValue "O36" "8.3 °C"
Set-TextValue "E37" "2026-02-10 07:49:50"
Set-TextValue "H37" "97%"
Set-TextValue "O37" "3.6 °C"
Set-TextValue "E38" "2026-02-10 07:49:52"
Set-TextValue "N38" "6.0 °C 7:05 TU"
Set-TextValue "O38" "7.4 °C"
Set-TextValue "E39" "2026-02-10 07:49:55"
Set-TextValue "K39" "0.1 MJ/m2"
Set-TextValue "M39" "1.2 °C 7:29 TU"
Set-TextValue "E40" "2026-02-10 07:49:57"
Set-TextValue "I40" "4.2 mm"
Set-TextValue "E41" "2026-02-10 07:50:00"
Set-TextValue "J41" "1005.2 hPa"
Set-TextValue "O41" "10.0 °C"
Set-TextValue "E42" "2026-02-10 07:50:02"
Set-TextValue "N42" "6.1 °C 7:16 TU"
Set-TextValue "O42" "7.7 °C"
Set-TextValue "E43" "2026-02-10 07:50:05"
Set-TextValue "N43" "4.6 °C 7:12 TU"
Set-TextValue "O43" "5.9 °C"
Set-TextValue "E44" "2026-02-10 07:50:07"
Set-TextValue "G44" "222 cm"
Set-TextValue "H44" "97%"
Set-TextValue "I44" "8.0 mm"
Set-TextValue "E45" "2026-02-10 07:50:09"
Set-TextValue "H45" "97%"
Set-TextValue "I45" "18.4 mm"
Set-TextValue "O45" "3.4 °C"
Set-TextValue "E46" "2026-02-10 07:50:12"
Set-TextValue "H46" "97%"

$ws.Range("ZZ1").ClearContents()
$excel.CutCopyMode = $false
